$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Esperado/Observado/valor p for rows 2-31 (unaffected by the row deletion below)
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0.27

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1

$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 0.06

$ws.Range("D6").Value = 59

$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.14

$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 1

$ws.Range("C11").Value = 38
$ws.Range("D11").Value = 27
$ws.Range("E11").Value = 0.01

$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 0.18

$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 0.15

$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0

$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 0.11

$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0.03

$ws.Range("D24").Value = 2

$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 13
$ws.Range("E27").Value = 0.01

$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0.27

$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 1

$ws.Range("C31").Value = 2
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0.14

# Remove obsolete row 32 (evento 720, Sindrome de rubeola congenita); rows below shift up one
$ws.Rows(32).Delete()

# Update Esperado/Observado/valor p for the shifted rows (now 32-38)
$ws.Range("D34").Value = 1

$ws.Range("C35").Value = 7
$ws.Range("D35").Value = 7
$ws.Range("E35").Value = 0.15

$ws.Range("C36").Value = 10
$ws.Range("D36").Value = 3

$ws.Range("C37").Value = 10
$ws.Range("E37").Value = 0.02
